# Rename sample identifiers on the "names" sheet and the "internal standard"
# sheet from the old Sx1d_phl_* / ZO4_phl_* labels to the new Sample_1..16
# labels, and update the remembered selections on those two sheets.

$wb = $excel.ActiveWorkbook

$sampleNames = @(
    "Sample_1", "Sample_2", "Sample_3", "Sample_4",
    "Sample_5", "Sample_6", "Sample_7", "Sample_8",
    "Sample_9", "Sample_10", "Sample_11", "Sample_12",
    "Sample_13", "Sample_14", "Sample_15", "Sample_16"
)

# --- "names" sheet: two blocks of 8 rows each (A3:A10, A13:A20) ---
$wsNames = $wb.Worksheets.Item("names")

for ($i = 0; $i -lt 8; $i++) {
    $row = 3 + $i
    $wsNames.Cells.Item($row, 1).Value = $sampleNames[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = 13 + $i
    $wsNames.Cells.Item($row, 1).Value = $sampleNames[$i + 8]
}

# --- "internal standard" sheet: two blocks of 8 rows (A3:A10, A12:A19) ---
$wsInternal = $wb.Worksheets.Item("internal standard")

for ($i = 0; $i -lt 8; $i++) {
    $row = 3 + $i
    $wsInternal.Cells.Item($row, 1).Value = $sampleNames[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $row = 12 + $i
    $wsInternal.Cells.Item($row, 1).Value = $sampleNames[$i + 8]
}

# --- update remembered selections ---
# "names" sheet selection moves to C28
[void]$wsNames.Range("C28").Select()

# restore "internal standard" as the active/selected sheet, then set its
# selection to A12:A19 (active cell A12)
[void]$wsInternal.Activate()
[void]$wsInternal.Range("A12:A19").Select()
